$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(13, 15, 16, 17, 18, 19, 20, 21, 22, 23, 26, 27)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 22).Value = 1
}

$ws.Columns.Item(22).ColumnWidth = 5.8

$ws.Range("V28").Select()
